$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7521266940440444
$ws.Range("C2").Value = 1.209530635601388
$ws.Range("D2").Value = 2.562311816749782
$ws.Range("E2").Value = 0.8176462137331231

$ws.Range("B3").Value = 0.8321548632773428
$ws.Range("C3").Value = 0.9953054322395707
$ws.Range("D3").Value = 2.380907242287568

$ws.Range("B4").Value = 0.6935836531676705
$ws.Range("C4").Value = 1.446089834638246
$ws.Range("D4").Value = 2.004859373446257
$ws.Range("E4").Value = 0.8127146885998382

$ws.Range("B5").Value = 0.3100022253933511
$ws.Range("E5").Value = 0.4933263677573559

$ws.Range("D6").Value = 1.378327106862104
$ws.Range("E6").Value = 0.5816833380721442

$ws.Range("B7").Value = 0.4291847764763447
$ws.Range("C7").Value = 0.1951495598821029
$ws.Range("D7").Value = 0.5655129046837359
$ws.Range("E7").Value = 0.8199105612194865

$ws.Range("B8").Value = 0.5765960342858791
$ws.Range("C8").Value = 1.98336826110382
$ws.Range("D8").Value = 4.70610678656399
$ws.Range("E8").Value = 0.6944421533041657

$ws.Range("B9").Value = 0.7314954403858797
$ws.Range("D9").Value = 4.180376060661897
$ws.Range("E9").Value = 0.7314954403858797

$ws.Range("B10").Value = 0.5171223429000441
$ws.Range("C10").Value = 1.785801060198129
$ws.Range("D10").Value = 3.958561777245468
$ws.Range("E10").Value = 0.6876838955551791

$ws.Range("D11").Value = 227.2050239264632
$ws.Range("E11").Value = 0.5486292772801611

$ws.Range("D12").Value = 200.7288507310735
$ws.Range("E12").Value = 0.6005781721361693

$ws.Range("B13").Value = 0.3286039516524467
$ws.Range("C13").Value = 80.57558914372228
$ws.Range("D13").Value = 177.6351139904597
$ws.Range("E13").Value = 0.6944014653713412

$ws.Range("D14").Value = 0.5735148163203111
$ws.Range("E14").Value = 0.453171177339487

$ws.Range("B15").Value = 0.5123697993192611
$ws.Range("D15").Value = 0.506951502460855
$ws.Range("E15").Value = 0.5123697993192606

$ws.Range("B16").Value = 0.3388455929761235
$ws.Range("C16").Value = 0.5353204702837886
$ws.Range("D16").Value = 0.03516888098594206
$ws.Range("E16").Value = 1.208968339218948
